# Update projection results on Sheet1 with refreshed simulation output values.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 (Year 1)
$ws.Range("B2").Value = 9478
$ws.Range("C2").Value = 9453
$ws.Range("D2").Value = 8370
$ws.Range("E2").Value = 0.8854331958108537
$ws.Range("F2").Value = 0.8830976999366955
$ws.Range("G2").Value = 0.09588380880107245
$ws.Range("H2").Value = 0.08467477101339696
$ws.Range("I2").Value = 40683573.30702042
$ws.Range("J2").Value = 14146954.60016621
$ws.Range("L2").Value = 14146954.60016621
$ws.Range("M2").Value = 54830527.90718663
$ws.Range("N2").Value = 800605047.7172
$ws.Range("O2").Value = 782905240.7132001
$ws.Range("P2").Value = 0.01767032900992073
$ws.Range("Q2").Value = 0.01806981721987046

# Row 3 (Year 2)
$ws.Range("B3").Value = 9666
$ws.Range("C3").Value = 9645
$ws.Range("D3").Value = 8566
$ws.Range("E3").Value = 0.8881285640228097
$ws.Range("F3").Value = 0.8861990482102214
$ws.Range("G3").Value = 0.09426000800759959
$ws.Range("H3").Value = 0.0835331293806226
$ws.Range("I3").Value = 42503835.53294694
$ws.Range("J3").Value = 14787168.70621267
$ws.Range("L3").Value = 14787168.70621267
$ws.Range("M3").Value = 57291004.23915961
$ws.Range("N3").Value = 837916296.6930281
$ws.Range("O3").Value = 820436120.6689579
$ws.Range("P3").Value = 0.01764754876420547
$ws.Range("Q3").Value = 0.01802354666461501

# Row 4 (Year 3)
$ws.Range("B4").Value = 9858
$ws.Range("C4").Value = 9832
$ws.Range("D4").Value = 8735
$ws.Range("E4").Value = 0.8884255492270138
$ws.Range("F4").Value = 0.886082369649016
$ws.Range("G4").Value = 0.09299641123507345
$ws.Range("H4").Value = 0.08240248043602824
$ws.Range("I4").Value = 44349172.10058596
$ws.Range("J4").Value = 15389914.85936836
$ws.Range("L4").Value = 15389914.85936836
$ws.Range("M4").Value = 59739086.95995432
$ws.Range("N4").Value = 874551910.723475
$ws.Range("O4").Value = 857102962.7175211
$ws.Range("P4").Value = 0.01759748583321603
$ws.Range("Q4").Value = 0.01795573639201207

# Row 5 (Year 4)
$ws.Range("B5").Value = 10054
$ws.Range("C5").Value = 10036
$ws.Range("D5").Value = 8932
$ws.Range("E5").Value = 0.889996014348346
$ws.Range("F5").Value = 0.8884026258205689
$ws.Range("G5").Value = 0.09156419815045283
$ws.Range("H5").Value = 0.08134587406801719
$ws.Range("I5").Value = 46337616.22484795
$ws.Range("J5").Value = 16051063.21331039
$ws.Range("L5").Value = 16051063.21331039
$ws.Range("M5").Value = 62388679.43815833
$ws.Range("N5").Value = 913766845.3762347
$ws.Range("O5").Value = 896280739.9132615
$ws.Range("P5").Value = 0.01756581921803205
$ws.Range("Q5").Value = 0.01790852184870529

# Row 6 (Year 5)
$ws.Range("B6").Value = 10254
$ws.Range("C6").Value = 10233
$ws.Range("D6").Value = 9093
$ws.Range("E6").Value = 0.8885957197302844
$ws.Range("F6").Value = 0.8867758923346987
$ws.Range("G6").Value = 0.09047387773774557
$ws.Range("H6").Value = 0.08023005366386975
$ws.Range("I6").Value = 48397322.24827884
$ws.Range("J6").Value = 16709871.56519883
$ws.Range("L6").Value = 16709871.56519883
$ws.Range("M6").Value = 65107193.81347767
$ws.Range("N6").Value = 954891245.1222031
$ws.Range("O6").Value = 937299419.2490215
$ws.Range("P6").Value = 0.01749924051619131
$ws.Range("Q6").Value = 0.01782767728436984
